$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (StatQuery) rows 2-5: unify on the corrected "samp.sample_site" count query
# (row 2 already had it; rows 3-5 had a stale "diag.primary_disease_site" version)
$cQuery = "MATCH (p:program)<--(s:study)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis), (c)<--(r:registration)`r`nMATCH (samp:sample)-->(c)`r`nWHERE s.clinical_study_designation IN ['MGT01'] and samp.sample_site in['Mammary Gland']`r`nOPTIONAL MATCH (cf:file)-[*]->(c)`r`nOPTIONAL MATCH (sf:file)-->(s)`r`nRETURN`r`n`tcount(distinct p) AS Programs,`r`n    count(distinct s) AS Studies,`r`n    count(distinct c) AS Cases,`r`n    count(distinct samp) AS Samples,`r`n    count(distinct cf) AS ``Case Files``,`r`n    count(distinct sf) AS ``Study Files```r`n`r`n    "
$ws.Range("C2").Value = $cQuery
$ws.Range("C3").Value = $cQuery
$ws.Range("C4").Value = $cQuery
$ws.Range("C5").Value = $cQuery

# Columns D & E (dbExcel / WebExcel filenames) rows 2-5: swap TC05 PrimaryDiseaseSite filenames
# for the TC44 SampleSite filenames
$neo4jFile = "TC44_Canine_Study_MGT01_SampleSite_MammaryGland_Neo4jData.xlsx"
$webFile = "TC44_Canine_Study_MGT01_SampleSite_MammaryGland_WebData.xlsx"
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile
$ws.Range("D5").Value = $neo4jFile
$ws.Range("E5").Value = $webFile

# Column B row 4 (FilesTab dbExcel query): add the sample match + switch filter to samp.sample_site
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`r`nWITH DISTINCT f, parent`r`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`r`nMATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`r`nMATCH (r:registration)-->(c)`r`nMATCH (f)-[*]->(samp:sample)`r`nWHERE s.clinical_study_designation IN ['MGT01'] and samp.sample_site in['Mammary Gland']`r`nOPTIONAL MATCH (f)-[*]->(samp:sample)`r`nWITH`r`n        DISTINCT f, parent, c, demo, diag, s, samp,`r`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`r`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`r`n        2 as precision`r`nWITH`r`n        f, parent, c, demo, diag, s, samp,`r`n        f.file_size /(1024^i) AS value,`r`n        10^precision AS factor,`r`n        units[i] as unit`r`nWITH`r`n        f, parent, c, demo, diag, s, samp, unit,`r`n        round(factor * value)/factor AS size`r`nRETURN`r`n        coalesce(f.file_name, '') AS ``File Name``,`r`n        coalesce(f.file_format, '') AS ``Format``,`r`n        coalesce(f.file_type, '') AS ``File Type``,`r`n        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`r`n        coalesce(labels(parent)[0], '') AS ``Association``,`r`n        coalesce(f.file_description, '') AS ``Description``,`r`n        coalesce(samp.sample_id, '') AS ``Sample ID``,`r`n        coalesce(c.case_id, '') AS ``Case ID``,`r`n        coalesce(demo.breed,'') AS Breed ,`r`n        coalesce(diag.disease_term,'') AS Diagnosis`r`n        order by f.file_name asc`r`n        limit 100"

# Move the active selection from D3 to D5
$ws.Range("D5").Select()
